# Fix locator for scenario 1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Order Info")

# --- Text / label updates (shared string changes) ---
$ws.Range("C14").Value = "sDS204-2311001"                 # Sales Order No.
$ws.Range("C19").Value = "27 Nov 2023 - 03 Dec 2023"      # Target Date range
$ws.Range("N22").Value = "04 Dec ~ 10 Dec"                # Forecast header week
$ws.Range("E24").Value = "pDS204-2311001"                 # Parts No. row 24
$ws.Range("E25").Value = "pDS204-2311001"                 # Parts No. row 25
$ws.Range("E26").Value = "pDS204-2311001"                 # Parts No. row 26

# --- Date updates ---
$ws.Range("C18").Value = 45250.0   # Order Date
$ws.Range("P23").Value = 45251.0   # Receiver Inbound Plan Date
$ws.Range("S23").Value = 45293.0
$ws.Range("T23").Value = 45331.0
$ws.Range("U23").Value = 45265.0   # Estimated Inbound Date

# --- Row 24 quantity shifts ---
$ws.Range("Q24").Value = 1620.0
$ws.Range("R24").Value = 0.0
$ws.Range("U24").Value = 1620.0

# --- Row 25 quantity shifts ---
$ws.Range("Q25").Value = 1620.0
$ws.Range("R25").Value = 0.0
$ws.Range("U25").Value = 1620.0

# --- Row 26 quantity shifts ---
$ws.Range("Q26").Value = 800.0
$ws.Range("R26").Value = 0.0
$ws.Range("U26").Value = 800.0
